$d = $word.ActiveDocument

# Update the simulation run timestamps/duration reported in the "Simulation parameters" section.

$d.Content.Find.Execute("Start time: 2017-12-27 18:32:59", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Start time: 2018-01-31 12:35:53", 2)

$d.Content.Find.Execute("End time: 2017-12-27 18:33:05", $true, $false, $false, $false, $false,
                         $true, 1, $false, "End time: 2018-01-31 12:36:01", 2)

$d.Content.Find.Execute("Duration: 5.84 secs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Duration: 8.72 secs", 2)
